# "update soal prioritas 1"
#
# 1. Remove the old "_GoBack" bookmark (it will be re-created on the new
#    title paragraph below); doing this first avoids ending up with two
#    bookmarks sharing the same name.
# 2. Insert a new centered title paragraph ("Soal Intoduction Flutter")
#    at the very top of the document, carrying the "_GoBack" bookmark.
# 3. Drop the stray <w:lastRenderedPageBreak/> that used to sit in front
#    of "Tampilkan".

$d = $word.ActiveDocument

# --- 1. remove the old "_GoBack" bookmark ----------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# --- 2. insert the new title paragraph -------------------------------------
$newParaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:ind w:left="360" w:hanging="360"/><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="32"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="32"/></w:rPr><w:t>Soal</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="32"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="32"/></w:rPr><w:t>Intoduction</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="32"/></w:rPr><w:t xml:space="preserve"> Flutter</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'

$startOfDoc = $d.Range(0, 0)
$startOfDoc.InsertXML($newParaXml)

# --- 3. strip the lastRenderedPageBreak marker before "Tampilkan" ---------
# A self-replace rewrites the run that holds the match, which drops the
# stray <w:lastRenderedPageBreak/> marker along with it.
$d.Content.Find.Execute("Tampilkan", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Tampilkan", 2) | Out-Null
